# Rename the "Text" sheet to "Sheet1" (example workbook renamed ahead of
# publishing the arb_excel package to pub.dev).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Text")
$ws.Name = "Sheet1"
